$wb = $excel.ActiveWorkbook

# Rename sheets (new randomized timestamp-based names)
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911495723593"
$wb.Worksheets.Item(2).Name = "NB_TO-16502911521125338"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911521135383"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911521770287"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911522652702"

# Sheet1: GNG_TO
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911495168178.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291149539339.csv"
$ws1.Range("B4").Value = "go_stims-16502911495413375.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911495713553.csv"

# Sheet2: NB_TO
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650291151628288.csv"
$ws2.Range("B3").Value = "TB-16502911520889387.csv"
$ws2.Range("B4").Value = "OB-16502911513389907.csv"
$ws2.Range("B5").Value = "OB-1650291151594819.csv"
$ws2.Range("B6").Value = "ZB-match_0-16502911498868396.csv"
$ws2.Range("B7").Value = "OB-16502911504969265.csv"
$ws2.Range("B8").Value = "ZB-match_1-16502911497485597.csv"
$ws2.Range("B9").Value = "TB-16502911516634002.csv"
$ws2.Range("B10").Value = "ZB-match_3-16502911501704128.csv"

# Sheet3: RS_TO - no cell value changes

# Sheet4: TOL_TO
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911521282945.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911521155338.csv"
$ws4.Range("B4").Value = "MM_stims-16502911521611154.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911521292925.csv"
$ws4.Range("B6").Value = "MM_stims-16502911521760652.csv"
$ws4.Range("B7").Value = "ZM_stims-1650291152162122.csv"

# Sheet5: vSAT_TO
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291152179026.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911522484515.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911522165935.csv"
$ws5.Range("B5").Value = "SAT_stims-1650291152192926.csv"
